$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '25.959.24'
$cell.Style = 'Normal'
$ws.Range('E2').Value = '  -1.39%  '
$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.639.59'
$cell.Style = 'Normal'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.42%  '
$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '215.56'
$cell.Style = 'Normal'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E8').Value = '  -0.75%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  -1.83%  '
$ws.Range('E11').Value = '  +0.12%  '
$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '1.866.06'
$cell.Style = 'Normal'
$ws.Range('E12').Value = '  -0.67%  '
$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '4.27'
$cell.Style = 'Normal'
$ws.Range('E13').Value = '  -0.85%  '
$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '1.658.00'
$cell.Style = 'Normal'
$ws.Range('E14').Value = '  -0.09%  '
$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '0.545'
$cell.Style = 'Normal'
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('E16').Value = '  -0.64%  '
$cell = $ws.Range('D17')
$cell.NumberFormat = '@'
$cell.Value = '63.02'
$cell.Style = 'Normal'
$ws.Range('E17').Value = '  -0.96%  '
$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '25.923.41'
$cell.Style = 'Normal'
$ws.Range('E18').Value = '  -1.49%  '
$cell = $ws.Range('D20')
$cell.NumberFormat = '@'
$cell.Value = '193.02'
$cell.Style = 'Normal'
$ws.Range('E20').Value = '  -2.04%  '
$ws.Range('E21').Value = '  -1.86%  '
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('E24').Value = '  +0.27%  '
$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '0.131'
$cell.Style = 'Normal'
$ws.Range('E25').Value = '  +4.36%  '
$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '143.50'
$cell.Style = 'Normal'
$ws.Range('E26').Value = '  +0.28%  '
$ws.Range('E27').Value = '  +0.37%  '
$ws.Range('E28').Value = '  -1.93%  '
$cell = $ws.Range('D29')
$cell.NumberFormat = '@'
$cell.Value = '15.60'
$cell.Style = 'Normal'
$ws.Range('E29').Value = '  -0.42%  '
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('E31').Value = '  -0.79%  '
$ws.Range('E32').Value = '  -1.60%  '
$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '3.25'
$cell.Style = 'Normal'
$ws.Range('E33').Value = '  -0.42%  '
$ws.Range('E34').Value = '  -5.10%  '
$ws.Range('E35').Value = '  +1.67%  '
$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '0.901'
$cell.Style = 'Normal'
$ws.Range('E36').Value = '  -1.81%  '
$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '1.135.48'
$cell.Style = 'Normal'
$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.544'
$cell.Style = 'Normal'
$ws.Range('E38').Value = '  -2.39%  '
$ws.Range('E39').Value = '  -1.32%  '
$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '0.0157'
$cell.Style = 'Normal'
$ws.Range('E40').Value = '  -0.35%  '
$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '5.49'
$cell.Style = 'Normal'
$ws.Range('E42').Value = '  -3.45%  '
$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '99.34'
$cell.Style = 'Normal'
$ws.Range('E43').Value = '  -1.17%  '
$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '0.797'
$cell.Style = 'Normal'
$ws.Range('E44').Value = '  -0.83%  '
$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '1.775.74'
$cell.Style = 'Normal'
$ws.Range('E45').Value = '  -0.68%  '
$ws.Range('E46').Value = '  +1.36%  '
$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '56.67'
$cell.Style = 'Normal'
$ws.Range('E47').Value = '  +0.30%  '
$ws.Range('E48').Value = '  +2.48%  '
$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '1.48'
$cell.Style = 'Normal'
$ws.Range('E49').Value = '  -2.13%  '
$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '7.70'
$cell.Style = 'Normal'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('E51').Value = '  -0.60%  '
